$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything currently in A:N (k:2..k:15
# headers in row 1, and the numeric fold data in rows 2-11) shifts right to B:O.
$ws.Range("A1").EntireColumn.Insert()

# The new column A needs "fold 1".."fold 10" labels in rows 2-11, styled like
# the row-1 headers (bold, centered, bordered). Copy the header cell's format
# so the new cells reuse the existing style record instead of minting a new one.
$ws.Range("B1").Copy()

$labels = @("fold 1","fold 2","fold 3","fold 4","fold 5","fold 6","fold 7","fold 8","fold 9","fold 10")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.PasteSpecial(-4122)
    $cell.Value = $labels[$i]
}

$excel.CutCopyMode = $false
